$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 658.125
$ws.Range("I5").Value = 658.125
$ws.Range("K5").Value = 658.125
$ws.Range("M5").Value = -543.125

$ws.Range("H17").Value = 334783
$ws.Range("J17").Value = 334783
$ws.Range("L17").Value = 1004349
$ws.Range("N17").Value = -1004685

$ws.Range("H18").Value = 1152.5
$ws.Range("I18").Value = 870
$ws.Range("K18").Value = 870
$ws.Range("M18").Value = -586

$ws.Range("H19").Value = 631.2778
$ws.Range("I19").Value = 468.2857
$ws.Range("K19").Value = 468.2857
$ws.Range("M19").Value = -293.2857

$ws.Range("H62").Value = 2669.8572
$ws.Range("I62").Value = 2669.8572
$ws.Range("K62").Value = 2669.8572
$ws.Range("M62").Value = -2045.8572

$ws.Range("H65").Value = 2669.8572
$ws.Range("I65").Value = 2669.8572
$ws.Range("K65").Value = 13349.286
$ws.Range("M65").Value = -10229.286

$ws.Range("H116").Value = 5014.4287
$ws.Range("I116").Value = 4850.25
$ws.Range("J116").Value = 5999.5
$ws.Range("K116").Value = 4850.25
$ws.Range("L116").Value = 5999.5
$ws.Range("M116").Value = -1408.25
$ws.Range("N116").Value = -12883.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 9199.25
$ws.Range("J13").Value = 17399.5
$ws.Range("L13").Value = 17399.5
$ws.Range("N13").Value = -17687.5

$ws.Range("H61").Value = 23814258
$ws.Range("I61").Value = 25645590
$ws.Range("J61").Value = 6948.6665
$ws.Range("K61").Value = 25645590
$ws.Range("L61").Value = 6948.6665
$ws.Range("M61").Value = -25645378
$ws.Range("N61").Value = -7372.6665

$ws.Range("H63").Value = 3896.5
$ws.Range("I63").Value = 3200
$ws.Range("K63").Value = 3200
$ws.Range("M63").Value = -2514

$ws.Range("H66").Value = 3896.5
$ws.Range("I66").Value = 3200
$ws.Range("K66").Value = 16000
$ws.Range("M66").Value = -12568

$ws.Range("H74").Value = 37041356
$ws.Range("I74").Value = 40004428
$ws.Range("K74").Value = 40004428
$ws.Range("M74").Value = -40003554

$ws.Range("H77").Value = 37041356
$ws.Range("I77").Value = 40004428
$ws.Range("K77").Value = 200022140
$ws.Range("M77").Value = -200017772

$ws.Range("H109").Value = 62500
$ws.Range("J109").Value = 62500
$ws.Range("L109").Value = 62500
$ws.Range("N109").Value = -65274

$ws.Range("H122").Value = 1681.1666
$ws.Range("I122").Value = 1276
$ws.Range("K122").Value = 3828
$ws.Range("M122").Value = -1378

$ws.Range("H136").Value = 23814258
$ws.Range("I136").Value = 25645590
$ws.Range("J136").Value = 6948.6665
$ws.Range("K136").Value = 76936770
$ws.Range("L136").Value = 20845.9995
$ws.Range("M136").Value = -76934220
$ws.Range("N136").Value = -25945.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 803.6667
$ws.Range("I94").Value = 803.6667
$ws.Range("K94").Value = 803.6667
$ws.Range("M94").Value = -352.6667

$ws.Range("H99").Value = 875.8570999999999
$ws.Range("I99").Value = 875.8570999999999
$ws.Range("K99").Value = 875.8570999999999
$ws.Range("M99").Value = 622.1429000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7786.1523
$ws.Range("I31").Value = 5325.6113
$ws.Range("J31").Value = 11637.435
$ws.Range("K31").Value = 5325.6113
$ws.Range("L31").Value = 11637.435
$ws.Range("M31").Value = -5030.6113
$ws.Range("N31").Value = -12227.435

$ws.Range("H34").Value = 7786.1523
$ws.Range("I34").Value = 5325.6113
$ws.Range("J34").Value = 11637.435
$ws.Range("K34").Value = 5325.6113
$ws.Range("L34").Value = 11637.435
$ws.Range("M34").Value = -5123.6113
$ws.Range("N34").Value = -12041.435

$ws.Range("H36").Value = 7024
$ws.Range("I36").Value = 4048
$ws.Range("J36").Value = 10000
$ws.Range("K36").Value = 4048
$ws.Range("L36").Value = 10000
$ws.Range("M36").Value = -3660
$ws.Range("N36").Value = -10776

$ws.Range("H40").Value = 7024
$ws.Range("I40").Value = 4048
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 4048
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -3888
$ws.Range("N40").Value = -10320

$ws.Range("H41").Value = 9346.916999999999
$ws.Range("I41").Value = 9346.916999999999
$ws.Range("K41").Value = 9346.916999999999
$ws.Range("M41").Value = -8918.916999999999

$ws.Range("H56").Value = 14000
$ws.Range("J56").Value = 14000
$ws.Range("L56").Value = 14000
# M56 intentionally left blank (no value)
$ws.Range("N56").Value = -15690

$ws.Range("H58").Value = 50013400
$ws.Range("I58").Value = 83350920
$ws.Range("K58").Value = 83350920
$ws.Range("M58").Value = -83350717

$ws.Range("H60").Value = 21574.715
$ws.Range("I60").Value = 5341.1113
$ws.Range("J60").Value = 33749.918
$ws.Range("K60").Value = 5341.1113
$ws.Range("L60").Value = 33749.918
$ws.Range("M60").Value = -4830.1113
$ws.Range("N60").Value = -34771.918

$ws.Range("H107").Value = 397010.03
$ws.Range("I107").Value = 543898.6
$ws.Range("J107").Value = 103232.8
$ws.Range("K107").Value = 543898.6
$ws.Range("L107").Value = 103232.8
$ws.Range("M107").Value = -541978.6
$ws.Range("N107").Value = -107072.8

$ws.Range("H136").Value = 50013400
$ws.Range("I136").Value = 83350920
$ws.Range("K136").Value = 250052760
$ws.Range("M136").Value = -250050210

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 125538.875
$ws.Range("J5").Value = 525
$ws.Range("L5").Value = 1575
$ws.Range("N5").Value = -1799

$ws.Range("H38").Value = 104
$ws.Range("I38").Value = 50
$ws.Range("J38").Value = 131
$ws.Range("K38").Value = 150
$ws.Range("L38").Value = 393
$ws.Range("M38").Value = 197
$ws.Range("N38").Value = -1087

$ws.Range("H135").Value = 125538.875
$ws.Range("J135").Value = 525
$ws.Range("L135").Value = 4725
$ws.Range("N135").Value = -9795

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2472.7334
$ws.Range("I80").Value = 2332.5833
$ws.Range("K80").Value = 2332.5833
$ws.Range("M80").Value = -1334.5833

$ws.Range("H83").Value = 2472.7334
$ws.Range("I83").Value = 2332.5833
$ws.Range("K83").Value = 11662.9165
$ws.Range("M83").Value = -6670.916499999999

$ws.Range("H92").Value = 6457.6665
$ws.Range("J92").Value = 6457.6665
$ws.Range("L92").Value = 6457.6665
$ws.Range("N92").Value = -10201.6665

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4083.3333
$ws.Range("I22").Value = 3375
$ws.Range("J22").Value = 5500
$ws.Range("K22").Value = 3375
$ws.Range("L22").Value = 5500
$ws.Range("M22").Value = -3080
$ws.Range("N22").Value = -6090

$ws.Range("H27").Value = 4083.3333
$ws.Range("I27").Value = 3375
$ws.Range("J27").Value = 5500
$ws.Range("K27").Value = 3375
$ws.Range("L27").Value = 5500
$ws.Range("M27").Value = -3268
$ws.Range("N27").Value = -5714

$ws.Range("H104").Value = 45666.332
$ws.Range("J104").Value = 45666.332
$ws.Range("L104").Value = 45666.332
$ws.Range("N104").Value = -52654.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2875.3333
$ws.Range("I81").Value = 2875.3333
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 5750.6666
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -4689.6666
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 2875.3333
$ws.Range("I84").Value = 2875.3333
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 28753.333
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -23449.333
$ws.Range("N84").ClearContents()
